$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain plain text,
# even when the string looks like a number (e.g. "0.9994"), and
# without leaving any number-format/style applied to the cell.
function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "26.898.21"
$ws.Range("E2").Value = "  -3.84%  "
Set-TextCell "D3" "1.731.55"
$ws.Range("E3").Value = "  -2.15%  "
Set-TextCell "D4" "0.9994"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextCell "D5" "309.75"
$ws.Range("E5").Value = "  -5.92%  "
Set-TextCell "D6" "0.9992"
$ws.Range("E6").Value = "  -0.12%  "
Set-TextCell "D7" "0.4927"
$ws.Range("E7").Value = "  +5.47%  "
Set-TextCell "D8" "0.3500"
$ws.Range("E8").Value = "  -0.71%  "
Set-TextCell "D9" "42.94"
$ws.Range("E9").Value = "  -2.12%  "
Set-TextCell "D10" "0.07252"
$ws.Range("E10").Value = "  -1.91%  "
Set-TextCell "D11" "1.054"
$ws.Range("E11").Value = "  -3.00%  "
Set-TextCell "D12" "0.9995"
$ws.Range("E12").Value = "  -0.07%  "
Set-TextCell "D13" "19.93"
$ws.Range("E13").Value = "  -3.52%  "
Set-TextCell "D14" "5.876"
$ws.Range("E14").Value = "  -2.41%  "
Set-TextCell "D15" "1.722.47"
$ws.Range("E15").Value = "  -2.51%  "
Set-TextCell "D16" "6.815"
$ws.Range("E16").Value = "  -5.38%  "
Set-TextCell "D17" "86.73"
$ws.Range("E17").Value = "  -5.97%  "
Set-TextCell "D18" "0.00001036"
$ws.Range("E18").Value = "  -2.01%  "
Set-TextCell "D19" "0.06398"
$ws.Range("E19").Value = "  -0.33%  "
Set-TextCell "D20" "0.9990"
$ws.Range("E20").Value = "  -0.10%  "
Set-TextCell "D21" "16.58"
$ws.Range("E21").Value = "  -2.15%  "
Set-TextCell "D22" "5.723"
$ws.Range("E22").Value = "  -1.43%  "
Set-TextCell "D23" "26.967.66"
$ws.Range("E23").Value = "  -3.74%  "
Set-TextCell "D24" "10.97"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("E25").Value = "  -4.74%  "
Set-TextCell "D26" "153.87"
$ws.Range("E26").Value = "  -6.74%  "
Set-TextCell "D27" "19.92"
$ws.Range("E27").Value = "  -0.78%  "
Set-TextCell "D28" "1.922.35"
$ws.Range("E28").Value = "  -2.42%  "
Set-TextCell "D29" "2.084"
$ws.Range("E29").Value = "  -5.46%  "
Set-TextCell "D30" "120.02"
$ws.Range("E30").Value = "  -2.70%  "
Set-TextCell "D31" "1.048"
$ws.Range("E31").Value = "  -2.54%  "
Set-TextCell "D32" "0.09335"
$ws.Range("E32").Value = "  -0.03%  "
Set-TextCell "D33" "3.578"
$ws.Range("E33").Value = "  -2.20%  "
Set-TextCell "D34" "5.392"
$ws.Range("E34").Value = "  -2.98%  "
Set-TextCell "D35" "0.05912"
$ws.Range("E35").Value = "  -3.29%  "
Set-TextCell "D36" "0.02183"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D37" "11.01"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D38" "1.428"
$ws.Range("E38").Value = "  -1.23%  "
Set-TextCell "D39" "4.749"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("E40").Value = "  -4.15%  "
Set-TextCell "D41" "0.9990"
$ws.Range("E41").Value = "  -0.13%  "
Set-TextCell "D42" "0.5999"
$ws.Range("E42").Value = "  -2.96%  "
Set-TextCell "D43" "1.111"
$ws.Range("E43").Value = "  -6.88%  "
Set-TextCell "D44" "7.439"
$ws.Range("E44").Value = "  -4.45%  "
Set-TextCell "D45" "12.84"
$ws.Range("E45").Value = "  -2.81%  "
Set-TextCell "D46" "3.574"
$ws.Range("E46").Value = "  -4.74%  "
Set-TextCell "D47" "0.5612"
$ws.Range("E47").Value = "  -3.41%  "
Set-TextCell "D48" "119.47"
$ws.Range("E48").Value = "  -3.72%  "
Set-TextCell "D49" "1.841"
$ws.Range("E49").Value = "  -4.95%  "
Set-TextCell "D50" "0.06648"
Set-TextCell "D51" "1.097"
$ws.Range("E51").Value = "  -2.74%  "
